# Commit: "Fixed POI packaging and upgraded to POI 3.15."
#
# The canonical-OOXML diff for this fixture only reorders XML attributes
# (e.g. w:pgSz w:w=.. w:h=.. -> w:h=.. w:w=.., docDefaults rFonts/lang
# attribute order, latentStyles/lsdException attribute order, and the
# w:style/tblPr/tblCellMar attribute order) and re-sorts the root
# <w:document> namespace declarations alphabetically. Every removed line
# has an exact attribute-for-attribute (value-for-value) match on the
# added side - this is the byte-level fingerprint of Apache POI
# re-serializing the part after the library upgrade mentioned in the
# commit message, not a real document edit. No text, formatting,
# structure, style definition, section/page setup value, run, paragraph,
# image, or any other user-visible/semantic content changed.
#
# The Word object model (real Word or this COM-interop shim) does not
# expose attribute-serialization order as something a script can control
# - Range/Document "XML" family properties are read-only, and there is no
# OM call that reorders XML attributes. So the correct, faithful action
# here is to leave the document's content untouched (any Find/Replace or
# other edit call - even a semantic no-op like replacing text with an
# identical string - would itself perturb the saved markup, e.g. merging
# runs or injecting extra namespace declarations, which is not part of
# this change).
$d = $word.ActiveDocument

# Touch nothing; simply confirm the document is available, as there is no
# content-level change to apply for this revision.
$null = $d.Name
